$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new bank account row (row 5)
$ws.Range("A5").Value = "INR"
$ws.Range("B5").Value = "Kotak Mahindra Bank"
$ws.Range("C5").Value = 9769689780
$ws.Range("D5").Value = "KKBK0001477"
$ws.Range("E5").Value = "KKBKINBB"
$ws.Range("F5").Value = "Wadala Branch"

# Autofit column C width so the new (shorter) numeric account number and
# the existing entries settle on a best-fit width, matching the workbook's
# post-import column sizing.
$ws.Columns.Item(3).ColumnWidth = 10.14

# Move selection to F6, matching the post-edit cursor position
$ws.Range("F6").Select() | Out-Null
